# Applies the transformation described by the diff:
# - Row 17 (old "El Señor de los Anillos III" entry, id U7sA0AEACAAJ) is removed.
# - Rows 18-20 shift up to become rows 17-19 (values unchanged, just moved up).
# - A brand-new row 20 is populated with "The Lord of the Rings Illustrated" (id ZcAlEAAAQBAJ).
#
# Net effect on cell values (rows 17..20, columns A..E):
#   Row17: UfYGAAAACAAJ | El señor de los anillos           | Desconocido           | 2002-02     | Kurt D. Bruner, Jim Ware
#   Row18: WmdWtQAACAAJ | El Señor de los anillos           | Desconocido           | 2002        | J. R. R. Tolkien
#   Row19: ZVwX0QEACAAJ | El Señor de los Anillos           | Desconocido           | 1985        | J. R. R. Tolkien
#   Row20: ZcAlEAAAQBAJ | The Lord of the Rings Illustrated | <long description>    | 2021-10-19  | J. R. R. Tolkien

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$longDescription = "A sumptuous slipcased edition of Tolkien's classic epic tale of adventure, fully illustrated in color by the author himself. This deluxe volume is quarterbound in leather and includes many special features unique to this edition. Since it was first published in 1954, The Lord of the Rings has been a book people have treasured. Steeped in unrivalled magic and otherworldliness, its sweeping fantasy and epic adventure has touched the hearts of young and old alike. Over 100 million copies of its many editions have been sold around the world, and occasional collectors' editions become prized and valuable items of publishing. This one-volume deluxe slipcased edition contains the complete text, fully corrected and reset, which is printed in red and black, and features thirty color illustrations, maps, and sketches drawn by Tolkien himself as he composed this epic work. These include the pages from the Book of Mazarbul, marvelous facsimiles created by Tolkien to accompany the famous `"Bridge of Khazad-dum`" chapter. Also appearing are two poster-size, fold-out maps revealing all the detail of Middle-earth. This very special deluxe edition is quarterbound in cloth and red leather, with raised ribs on the spine and stamped in two foils. The pages are edged in gold and contained within are special features unique to this edition."

# The "publishedDate" column (D) holds plain text values in the source file
# (e.g. "2002", "1985", "2021-10-19") rather than real numbers/dates. Excel
# would normally auto-convert such number- or date-looking text typed into
# a General formatted cell into an actual number/date. To keep these cells
# as text (matching the rest of the workbook), format the range as Text
# before assigning, then restore the cell style back to Normal afterwards
# so no visible formatting change is left behind.
$dateRange = $ws.Range("D17:D20")
$dateRange.NumberFormat = "@"

# Row 17
$ws.Cells.Item(17, 1).Value = "UfYGAAAACAAJ"
$ws.Cells.Item(17, 2).Value = "El señor de los anillos"
$ws.Cells.Item(17, 3).Value = "Desconocido"
$ws.Cells.Item(17, 4).Value = "2002-02"
$ws.Cells.Item(17, 5).Value = "Kurt D. Bruner, Jim Ware"

# Row 18
$ws.Cells.Item(18, 1).Value = "WmdWtQAACAAJ"
$ws.Cells.Item(18, 2).Value = "El Señor de los anillos"
$ws.Cells.Item(18, 3).Value = "Desconocido"
$ws.Cells.Item(18, 4).Value = "2002"
$ws.Cells.Item(18, 5).Value = "J. R. R. Tolkien"

# Row 19
$ws.Cells.Item(19, 1).Value = "ZVwX0QEACAAJ"
$ws.Cells.Item(19, 2).Value = "El Señor de los Anillos"
$ws.Cells.Item(19, 3).Value = "Desconocido"
$ws.Cells.Item(19, 4).Value = "1985"
$ws.Cells.Item(19, 5).Value = "J. R. R. Tolkien"

# Row 20
$ws.Cells.Item(20, 1).Value = "ZcAlEAAAQBAJ"
$ws.Cells.Item(20, 2).Value = "The Lord of the Rings Illustrated"
$ws.Cells.Item(20, 3).Value = $longDescription
$ws.Cells.Item(20, 4).Value = "2021-10-19"
$ws.Cells.Item(20, 5).Value = "J. R. R. Tolkien"

# Restore the default cell style on the date column so the workbook's
# formatting stays identical to the original (only the values changed).
$dateRange.Style = "Normal"
